$d = $word.ActiveDocument

# 1. Header contact line: add phone number, add http:// prefix to website
$d.Content.Find.Execute(
    "Paul Shorey  |  ps@artspaces.net  |  paulshorey.com",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Paul Shorey  |  ps@artspaces.net  |  321.446.5290  |  http://paulshorey.com",
    2) | Out-Null

# 2. Subtitle line: append "CSS ninja"
$d.Content.Find.Execute(
    "UI/UX designer  |  web developer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "UI/UX designer  |  web developer  |  CSS ninja",
    2) | Out-Null

# 3. Split the "Skills:" run into "Skills" + ":" runs (identical formatting, just a run break)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Skills:`r") {
        $full = $p.Range
        $colonRange = $d.Range($full.Start + 6, $full.Start + 7)
        $colonRange.Bold = 1
        $colonRange.Bold = 0
        break
    }
}

# 4. Design skills line rewording (scoped to its own paragraph)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Design - *") {
        $p.Range.Find.Execute(
            "Adobe CS, Sketch, SVG icons and fonts (custom and web services), responsive, video, new media",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "Adobe CS, Sketch, custom SVG icons and fonts, responsive, mobile, video, new media",
            2) | Out-Null
        break
    }
}

# 5/6. Development line: reword the text after "Development - " and keep the leading
#      space as its own run (matches the source run layout).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Development -*") {
        # Replace the long description first (still attached after " pr" at this point)
        $p.Range.Find.Execute(
            "efer and excel at anything to do with DOM / HTML / JSON, including CSS / SASS / GRUNT, have much experience in Javascript / JQuery / AngularJS / NodeJS, creating custom plugins and modules",
            $true, $false, $false, $false, $false, $true, 1, $false,
            "CSS3, SCSS, GRUNT, Javascript, JQuery, AngularJS, NodeJS, JSON, CORS, HTML5",
            2) | Out-Null

        # Collapse " pr" down to a single space
        $p.Range.Find.Execute(
            " pr",
            $true, $false, $false, $false, $false, $true, 1, $false,
            " ",
            2) | Out-Null

        # Force the trailing text back into its own run, separate from the leading space
        $full = $p.Range
        $splitStart = $full.Start + ("Development - ").Length
        $afterRange = $d.Range($splitStart, $full.End - 1)
        $afterRange.Bold = 1
        $afterRange.Bold = 0
        break
    }
}

# 7. Remove the trailing "Looking for a position..." paragraph entirely
$last = $d.Paragraphs($d.Paragraphs.Count)
if ($last.Range.Text -like "Looking for a position*") {
    $last.Range.Delete()
}
